$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3 data (staff member "Maryam"). In the source workbook every cell
# is stored as text (inline/shared string), even the numeric-looking ones
# (phone, age, salary), so force those columns to Text *before* assigning
# the values to stop Excel from auto-converting them to numbers/dates.
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 9).NumberFormat = "@"

$ws.Cells.Item(3, 1).Value = "maryam"
$ws.Cells.Item(3, 2).Value = "Maryam"
$ws.Cells.Item(3, 3).Value = "Maryam"
$ws.Cells.Item(3, 4).Value = "923432928333"
$ws.Cells.Item(3, 5).Value = "Islamabad"
$ws.Cells.Item(3, 6).Value = "2025-02-25"
$ws.Cells.Item(3, 7).Value = "20"
$ws.Cells.Item(3, 8).Value = "Female"
$ws.Cells.Item(3, 9).Value = "250000"
$ws.Cells.Item(3, 10).Value = "x"
$ws.Cells.Item(3, 11).Value = "x"
$ws.Cells.Item(3, 12).Value = "members,member_attendance,staff_attendance,payments,reports,staff,sales,inventory,packages"
$ws.Cells.Item(3, 13).Value = "manager"

# Drop the temporary "Text" number format again so the new row ends up with
# no explicit cell style, matching the rest of the (unstyled) data rows.
$ws.Range("A3:M3").ClearFormats()
